# Update "想去人数" (interest count) values in column F for both the
# "展览" (Worksheets item 1) and "全部类型" (Worksheets item 4) sheets,
# which carry duplicate copies of the same exhibition listing.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 11307
    4  = 10600
    7  = 756
    8  = 108
    12 = 10511
    20 = 11079
    21 = 10824
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
